$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.666.15"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "1.896.60"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5264"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3795"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07234"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9007"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.06%  "

$ws.Range("D12").Value = "1.913.22"
$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07625"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.429"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9993"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008665"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9995"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "27.702.52"
$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.143"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").Value = "2.166.95"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.859"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.166"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("E30").Value = "  -2.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.813"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09150"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05263"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.158"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.220"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7718"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02079"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.567"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.074"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.29%  "

$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5555"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.682"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.713"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1509"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4793"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("E47").Value = "  -2.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9991"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.585"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.41%  "

$ws.Range("E50").Value = "  -1.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "37.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
